# Hortaliza, Femacal de La Calera - Ciboulette
# Insert a new weekly record at the top of the data block (row 99),
# pushing every existing row down by one (the last row's data falls
# into the newly created row 249).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 99:248 down to 100:249, leaving a blank row 99 behind.
$ws.Rows("99:99").Insert()

# Populate the new row 99 with the latest observation.
$ws.Range("A99").Value2 = 3
$ws.Range("B99").Value2 = "Femacal de La Calera"
$ws.Range("C99").Value2 = "Coquimbo"
$ws.Range("D99").Value2 = 44580
$ws.Range("E99").Value2 = 5
$ws.Range("F99").Value2 = 100112039
$ws.Range("G99").Value2 = "Ciboulette"
$ws.Range("H99").Value2 = "Sin especificar"
$ws.Range("I99").Value2 = "Primera"
$ws.Range("J99").Value2 = 160
$ws.Range("K99").Value2 = 1500
$ws.Range("L99").Value2 = 1500
$ws.Range("M99").Value2 = 1500
$ws.Range("N99").Value2 = "`$/docena de atados"
$ws.Range("O99").Value2 = "Provincia de Quillota"
$ws.Range("P99").Value2 = 500
$ws.Range("Q99").Value2 = 3
$ws.Range("R99").Value2 = "Hortaliza"
